$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("I2").Value = 1.86
$ws.Range("Q2").Value = 1.68

$ws.Range("I3").Value = 3.75

$ws.Range("F4").Value = 1.6
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 7.4
$ws.Range("K4").Value = 4.6
$ws.Range("P4").Value = 1.99
$ws.Range("Q4").Value = 1.81

$ws.Range("F5").Value = 1.78
$ws.Range("G5").Value = 1.99
$ws.Range("H5").Value = 4.6
$ws.Range("I5").Value = 6.2
$ws.Range("J5").Value = 3.3
$ws.Range("M5").Value = 1.07
$ws.Range("P5").Value = 1.79
$ws.Range("Q5").Value = 1.92
$ws.Range("AG5").Value = 970

$ws.Range("F6").Value = 1.21
$ws.Range("G6").Value = 1.29
$ws.Range("H6").Value = 14.5
$ws.Range("I6").Value = 20
$ws.Range("K6").Value = 8.800000000000001
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 1.35

$ws.Range("F7").Value = 5.6
$ws.Range("H7").Value = 1.17
$ws.Range("I7").Value = 1.21
$ws.Range("J7").Value = 8
$ws.Range("K7").Value = 11.5
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.35

$ws.Range("F8").Value = 2.66
$ws.Range("G8").Value = 3.15
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.95
$ws.Range("P8").Value = 1.38
$ws.Range("Q8").Value = 2.84

$ws.Range("F9").Value = 1.81
$ws.Range("G9").Value = 1.93
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 2.88
$ws.Range("K9").Value = 3.9
$ws.Range("P9").Value = 1.58
$ws.Range("Q9").Value = 2.44

$ws.Range("F10").Value = 1.77
$ws.Range("G10").Value = 2.04
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 6.6
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 3.9
$ws.Range("P10").Value = 1.68
$ws.Range("Q10").Value = 2.16

$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 3.85

$ws.Range("F14").Value = 1.68
$ws.Range("G14").Value = 1.7
$ws.Range("H14").Value = 5.9
$ws.Range("I14").Value = 6.4
$ws.Range("J14").Value = 3.9
